$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 265, shifting the existing rows
# (old 265..377) down to (267..379).
$ws.Rows.Item(265).Insert()
$ws.Rows.Item(265).Insert()

# Row 265: new "Primera" quality entry for Betarraga @ Feria Lagunitas de Puerto Montt
$ws.Cells.Item(265, 1).Value = 4
$ws.Cells.Item(265, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(265, 3).Value = "Los Lagos"
$ws.Cells.Item(265, 4).Value = 44875
$ws.Cells.Item(265, 5).Value = 10
$ws.Cells.Item(265, 6).Value = 100114014
$ws.Cells.Item(265, 7).Value = "Betarraga"
$ws.Cells.Item(265, 8).Value = "Sin especificar"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 250
$ws.Cells.Item(265, 11).Value = 1500
$ws.Cells.Item(265, 12).Value = 1500
$ws.Cells.Item(265, 13).Value = 1500
$ws.Cells.Item(265, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(265, 15).Value = "Región del Maule"
$ws.Cells.Item(265, 16).Value = 300
$ws.Cells.Item(265, 17).Value = 5
$ws.Cells.Item(265, 18).Value = "Hortaliza"

# Row 266: new "Segunda" quality entry for Betarraga @ Feria Lagunitas de Puerto Montt
$ws.Cells.Item(266, 1).Value = 4
$ws.Cells.Item(266, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(266, 3).Value = "Los Lagos"
$ws.Cells.Item(266, 4).Value = 44875
$ws.Cells.Item(266, 5).Value = 10
$ws.Cells.Item(266, 6).Value = 100114014
$ws.Cells.Item(266, 7).Value = "Betarraga"
$ws.Cells.Item(266, 8).Value = "Sin especificar"
$ws.Cells.Item(266, 9).Value = "Segunda"
$ws.Cells.Item(266, 10).Value = 250
$ws.Cells.Item(266, 11).Value = 1200
$ws.Cells.Item(266, 12).Value = 1200
$ws.Cells.Item(266, 13).Value = 1200
$ws.Cells.Item(266, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(266, 15).Value = "Región del Maule"
$ws.Cells.Item(266, 16).Value = 240
$ws.Cells.Item(266, 17).Value = 5
$ws.Cells.Item(266, 18).Value = "Hortaliza"
